# This script reshuffles the weekly price-report rows (rows 2-13) of the
# "Hortaliza, Mapocho Venta Directa de Santiago - Zapallo italiano" sheet.
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) are permuted across
# rows 2-13 (row 3 keeps its original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together as one "record" when rows are reshuffled.
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values for every source row so that
# later writes don't clobber data we still need to read.
$snapshot = @{}
foreach ($row in 2..13) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Mapping of target row -> source row (row 3 is unchanged, so it is omitted).
$mapping = @{
    2  = 9
    4  = 6
    5  = 12
    6  = 8
    7  = 13
    8  = 7
    9  = 5
    10 = 4
    11 = 10
    12 = 2
    13 = 11
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $sourceData[$col]
    }
}
